# Add WEEK 4 task list to the Tasks table (Table1) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Table1")

# New WEEK 4 rows: [Task, Owner]; Week = "WEEK 4", Status = "Pending" for all.
$newTasks = @(
    @("Add New Machine Learning Techniques: LightGBM, CatBoost, Random forest With Random and Grid search for all of these", "Abhik"),
    @("SVC try with different kernel in grid search", "Abhik"),
    @("In the accuracy score table also add scores from data without augmentation", "Abhik"),
    @("Model interpretation - Lime and other packages", "Abhik"),
    @("Add Embedding projection using Tensor board", "Abhijit"),
    @("Clustering of data and visualization: Topic Modeling (use LDAPyviz), Document clustering", "Abhijit"),
    @("Add visualization for train, test data and show all groups are present in both set", "Abhijit"),
    @("Deep learning hypeparameter tuning using packages", "Abhijit")
)

$firstNewRow = $true
foreach ($task in $newTasks) {
    $row = $lo.ListRows.Add()
    $r = $row.Range.Row

    $ws.Cells.Item($r, 1).Value = "WEEK 4"
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 2).Value = $task[0]
    if ($firstNewRow) {
        $ws.Cells.Item($r, 2).WrapText = $true
        $ws.Rows.Item($r).RowHeight = 30
        $firstNewRow = $false
    }

    $ws.Cells.Item($r, 3).Value = $task[1]
    $ws.Cells.Item($r, 4).Value = "Pending"
}

$ws.Range("D31").Select() | Out-Null
